$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as row 21 (pushing the
# previous row 21 down to row 22, preserving its data/formatting untouched).
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new weekly record.
$ws.Cells.Item(21, 1).Value = 6
$ws.Cells.Item(21, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(21, 3).Value = "Metropolitana"
$ws.Cells.Item(21, 4).Value = 44505
$ws.Cells.Item(21, 5).Value = 13
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100102
$ws.Cells.Item(21, 8).Value = "Cítricos"
$ws.Cells.Item(21, 9).Value = 100102006
$ws.Cells.Item(21, 10).Value = "Pomelo"
$ws.Cells.Item(21, 11).Value = "Start Ruby"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 15
$ws.Cells.Item(21, 14).Value = 150000
$ws.Cells.Item(21, 15).Value = 150000
$ws.Cells.Item(21, 16).Value = 150000
$ws.Cells.Item(21, 17).Value = "`$/bins (350 kilos)"
$ws.Cells.Item(21, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(21, 19).Value = 429
$ws.Cells.Item(21, 20).Value = 350
